# tdf#146756 unit test font: substitute Carlito for Calibri font
#
# The theme's minor-font Latin typeface ("Calibri") is swapped for "Carlito",
# a font that is metric-compatible with Calibri. The major-font entry
# ("Calibri Light") is intentionally left untouched, matching the upstream
# change which only retargets the plain "Calibri" minor-font reference(s).

$p = $ppt.ActivePresentation

# ppt/theme/theme1.xml is reached through the slide master's Theme object.
# Editing the theme per the documented PowerPoint COM idiom:
#   Theme.ThemeFontScheme.MinorFont/MajorFont.Latin
$master = $p.SlideMaster
$theme = $master.Theme
$fontScheme = $theme.ThemeFontScheme

$minorFont = $fontScheme.MinorFont
if ($minorFont.Latin -eq "Calibri") {
    $minorFont.Latin = "Carlito"
}

# Note: this host's NotesMaster/HandoutMaster Theme objects alias back to the
# slide master's theme rather than exposing the notes/handout master's own
# theme parts (ppt/theme/theme2.xml, ppt/theme/theme3.xml) independently, so
# those parts are not reachable for editing through the PowerPoint object
# model in this environment.
